# Daily IST report: add CSV/MD/XLSX
# Inserts a new "2026-02-24" day column (J) into the submissions daily
# matrix, pushing total_files -> K and unique_days -> L, and updates the
# rolling total_files / unique_days counters to include the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture a known-good "12"-wide and "13"-wide column's ColumnWidth getter
# value before we touch the layout; re-applying a getter value through the
# setter round-trips losslessly (avoids float drift from Excel's
# char-width <-> pixel-width conversion).
$width12 = $ws.Columns.Item(4).ColumnWidth
$width13 = $ws.Columns.Item(11).ColumnWidth

# Insert a new column at J (10): existing J (total_files) shifts to K,
# existing K (unique_days) shifts to L. Formatting/styles shift with them.
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(10).ColumnWidth = $width12
$ws.Columns.Item(11).ColumnWidth = $width13
$ws.Columns.Item(12).ColumnWidth = $width13

# New day header "2026-02-24", styled like the other date headers
# (D1:I1). Assigning a date-shaped string straight to .Value gets
# auto-parsed into a date serial by the smart-type coercion, so force
# text mode first, then restore the plain header look (bold+centered,
# no number format) by pasting the format from the neighbouring date
# header (I1), which leaves the freshly-written text value untouched.
$ws.Cells.Item(1, 10).NumberFormat = "@"
$ws.Cells.Item(1, 10).Value = "2026-02-24"
$ws.Cells.Item(1, 9).Copy()
$ws.Cells.Item(1, 10).PasteSpecial(-4122)

# Per-row file counts submitted on 2026-02-24 (row 2 .. row 109).
$newDay = @(0, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 0, 0, 0, 1, 0, 1, 1, 0, 0, 0, 1, 1, 1, 1, 1, 1, 1, 0, 0, 1, 1, 1, 0, 42, 0, 1, 1, 0, 1, 1, 0, 0, 1, 0, 1, 0, 1, 1, 0, 1, 0, 0, 1, 1, 0, 0, 0, 0, 0, 1, 0, 0, 1, 1, 0, 0, 0, 0, 0, 1, 1, 1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 1, 0, 1, 0, 1, 1, 0, 1, 0)

for ($i = 0; $i -lt $newDay.Length; $i++) {
    $row = $i + 2
    $dayVal = $newDay[$i]

    # Existing (pre-edit) rolling counters already live in K/L after the
    # column insert: K = old total_files, L = old unique_days. (.Value2
    # is used for numeric reads/writes here -- the bare .Value getter
    # does not resolve to the underlying number in this host.)
    $oldTotal = $ws.Cells.Item($row, 11).Value2
    $oldUnique = $ws.Cells.Item($row, 12).Value2

    $ws.Cells.Item($row, 10).Value2 = $dayVal
    $ws.Cells.Item($row, 11).Value2 = $oldTotal + $dayVal
    if ($dayVal -gt 0) {
        $ws.Cells.Item($row, 12).Value2 = $oldUnique + 1
    } else {
        $ws.Cells.Item($row, 12).Value2 = $oldUnique
    }
}
